# NATMI TPM update: Vegfc-Kdr.xlsx
# The underlying ligand/receptor TPM inputs were recomputed upstream, which
# changes the expressing-cell counts, detection rates, average/total
# expression values and all of the derived specificity / edge-weight
# columns (E:T, excluding the untouched K/L "Receptor-expressing
# cells"/"Receptor detection rate" columns) for every sending/target
# cluster combination on the sheet. Apply the refreshed values cell by cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.055887666666667
$ws.Range("H2").Value = 12.167663
$ws.Range("I2").Value = 0.4763357569530485
$ws.Range("J2").Value = 0.4763357569530485
$ws.Range("M2").Value = 211.980367
$ws.Range("N2").Value = 635.9411009999999
$ws.Range("O2").Value = 0.9885149156420702
$ws.Range("P2").Value = 0.9885149156420702
$ws.Range("Q2").Value = 859.7685560907736
$ws.Range("R2").Value = 7737.917004816963
$ws.Range("S2").Value = 0.4708650006017444
$ws.Range("T2").Value = 0.4708650006017444

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.055887666666667
$ws.Range("H3").Value = 12.167663
$ws.Range("I3").Value = 0.4763357569530485
$ws.Range("J3").Value = 0.4763357569530485
$ws.Range("O3").Value = 0.003992992409159323
$ws.Range("P3").Value = 0.003992992409159324
$ws.Range("Q3").Value = 3.472936284299223
$ws.Range("R3").Value = 31.256426558693
$ws.Range("S3").Value = 0.001902005061724683
$ws.Range("T3").Value = 0.001902005061724683

# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.055887666666667
$ws.Range("H4").Value = 12.167663
$ws.Range("I4").Value = 0.4763357569530485
$ws.Range("J4").Value = 0.4763357569530485
$ws.Range("O4").Value = 0.007492091948770576
$ws.Range("P4").Value = 0.007492091948770576
$ws.Range("Q4").Value = 6.516305394046445
$ws.Range("R4").Value = 58.64674854641801
$ws.Range("S4").Value = 0.003568751289579473
$ws.Range("T4").Value = 0.003568751289579473

# Row 5: FAPs -> ECs
$ws.Range("I5").Value = 0.4564376967244237
$ws.Range("J5").Value = 0.4564376967244237
$ws.Range("M5").Value = 211.980367
$ws.Range("N5").Value = 635.9411009999999
$ws.Range("O5").Value = 0.9885149156420702
$ws.Range("P5").Value = 0.9885149156420702
$ws.Range("Q5").Value = 823.8532877909422
$ws.Range("R5").Value = 7414.67959011848
$ws.Range("S5").Value = 0.4511954712734045
$ws.Range("T5").Value = 0.4511954712734045

# Row 6: FAPs -> FAPs
$ws.Range("I6").Value = 0.4564376967244237
$ws.Range("J6").Value = 0.4564376967244237
$ws.Range("O6").Value = 0.003992992409159323
$ws.Range("P6").Value = 0.003992992409159324
$ws.Range("S6").Value = 0.001822552258274789
$ws.Range("T6").Value = 0.001822552258274789

# Row 7: FAPs -> MuSCs
$ws.Range("I7").Value = 0.4564376967244237
$ws.Range("J7").Value = 0.4564376967244237
$ws.Range("O7").Value = 0.007492091948770576
$ws.Range("P7").Value = 0.007492091948770576
$ws.Range("S7").Value = 0.003419673192744441
$ws.Range("T7").Value = 0.003419673192744441

# Row 8: MuSCs -> ECs
$ws.Range("G8").Value = 0.5724183333333334
$ws.Range("I8").Value = 0.06722654632252778
$ws.Range("J8").Value = 0.06722654632252777
$ws.Range("M8").Value = 211.980367
$ws.Range("N8").Value = 635.9411009999999
$ws.Range("O8").Value = 0.9885149156420702
$ws.Range("P8").Value = 0.9885149156420702
$ws.Range("Q8").Value = 121.3414483775283
$ws.Range("R8").Value = 1092.073035397755
$ws.Range("S8").Value = 0.06645444376692126
$ws.Range("T8").Value = 0.06645444376692126

# Row 9: MuSCs -> FAPs
$ws.Range("G9").Value = 0.5724183333333334
$ws.Range("I9").Value = 0.06722654632252778
$ws.Range("J9").Value = 0.06722654632252777
$ws.Range("O9").Value = 0.003992992409159323
$ws.Range("P9").Value = 0.003992992409159324
$ws.Range("Q9").Value = 0.4901448370894445
$ws.Range("S9").Value = 0.000268435089159851
$ws.Range("T9").Value = 0.000268435089159851

# Row 10: MuSCs -> MuSCs
$ws.Range("G10").Value = 0.5724183333333334
$ws.Range("I10").Value = 0.06722654632252778
$ws.Range("J10").Value = 0.06722654632252777
$ws.Range("O10").Value = 0.007492091948770576
$ws.Range("P10").Value = 0.007492091948770576
$ws.Range("Q10").Value = 0.9196637036588889
$ws.Range("S10").Value = 0.0005036674664466626
$ws.Range("T10").Value = 0.0005036674664466625
